$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: quantity/value updated, client cell (D2) re-pointed to new OCR match "opini" with default style, E2 cleared ---
$ws.Range("A2").Value = 5
$ws.Range("C2").Value = 19
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "opini"
$ws.Range("E2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 6

# --- Row 4 ---
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 19

# --- Row 5: cleared out (kept styles) ---
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()

# --- Row 6: cleared out (kept styles) ---
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# --- New OCR lookup table columns I/J for rows 13-18 (lowercase variant in I, uppercase in J) ---
$ws.Range("J13").Value = "BELINELLI 2"
$ws.Range("J15").Value = "BELINELLI 3"
$ws.Range("J14").Value = "BELINELLI 1"
$ws.Range("I13").Value = "belinelli 2"
$ws.Range("I14").Value = "belinelli 1"
$ws.Range("I15").Value = "belinelli 3"
$ws.Range("I16").Value = "belinelli 4"
$ws.Range("J16").Value = "BELINELLI 4"
$ws.Range("J17").Value = "OPINI OPINI"
$ws.Range("J18").Value = "TITI"
$ws.Range("I18").Value = "titi"
$ws.Range("I17").Value = "opini opini"

# --- Column widths (best-effort autofit to mirror the widened I/J OCR columns + resized data columns) ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Range("I1:J1").Columns.AutoFit()

# --- View state: scroll so column B is leftmost, select the last-entered OCR cell ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$ws.Range("I18").Select()
